$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column B (RPs) values to reflect the new normalization method
$ws.Range("B3").Value = 0.7789538212688253
$ws.Range("B4").Value = 7.931993287521896
$ws.Range("B5").Value = 14.800122604107681
$ws.Range("B6").Value = 15.863986575043793
$ws.Range("B7").Value = 23.79597986256569
$ws.Range("B8").Value = 28.042337565677713
$ws.Range("B9").Value = 31.727973150087585
$ws.Range("B10").Value = 39.65996643760948
